$d = $word.ActiveDocument

# Locate the "ZEP" Heading2 paragraph, then remove the paragraph that
# immediately follows it (the italic "Sefania" subtitle paragraph),
# merging what's left of the following paragraph's tail run(s) with it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ZEP`r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $sefaniaPara = $target.Next()
    if ($sefaniaPara.Range.Text -eq "Sefania`r") {
        $sefaniaPara.Range.Delete()
    }
}
